# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table -------------------------------------------------
# The "22.240.0.6" driver (old row 4) dropped out of the bad-driver list
# this week; delete its row and let everything below shift up.
$ws.Rows("4:4").Delete()

# Remaining bad-driver rows, reordered and refreshed with this week's figures.
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.160.0.4"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 94.90000000000001

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 218
$ws.Range("D4").Value = 97.5

# Totals row (now row 5 after the deletion above)
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 220

# --- Good Drivers table (now starts at row 11 after the deletion) ------
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B13").Value = 56018
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = 0

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = 0

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B15").Value = 442178
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "2024-11-10"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B16").Value = 77849
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "2021-08-18"

# Row 17 (21.110.3.2) is unchanged this week.

$ws.Range("B18").Value = 113652
$ws.Range("E18").Value = "2019-12-14"
